$wb = $excel.ActiveWorkbook

# Sheet 2 (Tabelle2): add header/data columns
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("A1").Value = "column1"
$ws2.Range("B1").Value = "column2"
$ws2.Range("A2").Value = "data1"
$ws2.Range("B2").Value = "data2"
$null = $ws2.Range("B2").Select()

# Sheet 3 (Tabelle3): add a row of data
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("A1").Value = "this"
$ws3.Range("B1").Value = "thing"
$ws3.Range("C1").Value = "is"
$ws3.Range("D1").Value = "legit"
$null = $ws3.Range("D1").Select()

# Make Tabelle3 the active sheet/tab
$null = $ws3.Activate()
